# Transition rule summary tables: add "Within 5 miles" and "Within 10 miles"
# of HFC production facility columns (F and G) to both the "Means" and
# "Standard Deviations" sheets, and refresh the previously-existing values
# that changed as part of the underlying script re-run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Means"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Means")

# New column headers
$ws1.Cells.Item(1, 6).Value = "Within 5 miles of HFC production facility"
$ws1.Cells.Item(1, 7).Value = "Within 10 miles of HFC production facility"

# Row labels (A2:A10) stay the same; only refresh the data columns B:G
$means = @(
    # rowIndex, B,    C,    D,    E,    F,    G
    @(2,  72,   60,   73,   65,   58,   49),
    @(3,  13,   5.8,  2.1,  3,    3.9,  3.6),
    @(4,  15,   35,   25,   32,   39,   47),
    @(5,  18,   39,   36,   44,   50,   55),
    @(6,  71,   83,   88,   88,   83,   80),
    @(7,  7.3,  7.3,  3.5,  4.8,  6,    6.5),
    @(8,  5.8,  5.8,  5.6,  4.1,  5,    4.6),
    @(9,  29,   31,   40,   38,   38,   37),
    @(10, 0.37, 0.43, 0.44, 0.45, 0.44, 0.46)
)

foreach ($row in $means) {
    $r = $row[0]
    for ($c = 1; $c -le 6; $c++) {
        $ws1.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# ---------------------------------------------------------------------
# Sheet 2: "Standard Deviations"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Standard Deviations")

# New column headers
$ws2.Cells.Item(1, 6).Value = "Within 5 mile of HFC production facility SD"
$ws2.Cells.Item(1, 7).Value = "Within 10 mile of HFC production facility SD"

$sds = @(
    # rowIndex, B,    C,    D,     E,    F,     G
    @(2,  27,   23,   12,    16,   17,    19),
    @(3,  23,   10,   1.5,   2.8,  4,     4),
    @(4,  16,   21,   12,    16,   16,    19),
    @(5,  22,   28,   16,    20,   24,    25),
    @(6,  37,   44,   19,    27,   30,    28),
    @(7,  8.7,  8.2,  4,     5.7,  5.9,   6.4),
    @(8,  7.8,  6.7,  13,    6.8,  6.2,   4.9),
    @(9,  10,   8.9,  0,     3.8,  4.1,   4.5),
    @(10, 0.14, 0.23, 0.051, 0.05, 0.052, 0.052)
)

foreach ($row in $sds) {
    $r = $row[0]
    for ($c = 1; $c -le 6; $c++) {
        $ws2.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
